# Applies the reordering of "step 3" content among TC3/TC4/TC5, and
# swaps the "step 2" Expected Result content between TC7/TC8, as
# described by the commit's XML diff (v1.2.1 -> v1.2.3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# --- TC3 / TC4 / TC5 block ---
# Before:
#   TC3 (row 32): "Chefe Clica em visualizar comprovante." / "SYSTEM Exibe modal com o comprovante."
#   TC4 (row 41): "Chefe Clica para detalhar a solicitação de diária." / "SYSTEM Apresenta a tela de Detalhar Diárias"
#   TC5 (row 50): "Chefe Clica em excluir comprovante." / "SYSTEM Exclui o comprovante."
# After (rotate so each label keeps its row, content shifts TC4->TC3, TC5->TC4, TC3->TC5):
#   TC3 (row 32): "Chefe Clica para detalhar a solicitação de diária." / "SYSTEM Apresenta a tela de Detalhar Diárias"
#   TC4 (row 41): "Chefe Clica em excluir comprovante." / "SYSTEM Exclui o comprovante."
#   TC5 (row 50): "Chefe Clica em visualizar comprovante." / "SYSTEM Exibe modal com o comprovante."

$ws.Range("B32").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

$ws.Range("B41").Value = "Chefe Clica em excluir comprovante."
$ws.Range("D41").Value = "SYSTEM Exclui o comprovante."

$ws.Range("B50").Value = "Chefe Clica em visualizar comprovante."
$ws.Range("D50").Value = "SYSTEM Exibe modal com o comprovante."

# --- TC7 / TC8 block ---
# Before:
#   TC7 (row 67): "SYSTEM Identifica que a solicitação indicada pelo usuário ainda não pode ter sua prestação de
#                  contas realizada; Exibe mensagem de erro (MSG212 - Prestação de contas ainda não pode ser
#                  realizada) para o usuário, impedindo que ele preste contas (anexa arquivos e etc)."
#   TC8 (row 75): "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses
#                  dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Permite não permite um novo envio ou
#                  alterações na prestação (exclusão de documentos)."
# After (swap the two Expected Result texts between TC7 and TC8):

$ws.Range("D67").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Permite não permite um novo envio ou alterações na prestação (exclusão de documentos)."
$ws.Range("D75").Value = "SYSTEM Identifica que a solicitação indicada pelo usuário ainda não pode ter sua prestação de contas realizada; Exibe mensagem de erro (MSG212 - Prestação de contas ainda não pode ser realizada) para o usuário, impedindo que ele preste contas (anexa arquivos e etc)."
